$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.848
$ws.Range("E3").Value = 12.974
$ws.Range("E5").Value = 13.165
$ws.Range("B9").Value = 6.795
$ws.Range("E11").Value = 13.012
$ws.Range("E12").Value = 13
$ws.Range("B13").Value = 6.472
$ws.Range("B16").Value = 5.825
$ws.Range("B18").Value = 6.313
$ws.Range("B20").Value = 6.661
$ws.Range("E21").Value = 12.988
